# Applies the "fixed giftcards" edit to the GiftCards sheet.
#
# Summary of the change (see commit message: "fixed giftcards, added
# CreateAndRetrieveBooking script, added global variables"):
#   - TCNumber in row 2 (A2) becomes a text value "1" instead of a number.
#   - The ContainsValidation / NotContainsValidation columns (H, I) for the
#     "Valid Card Number" (row 2) and "Invalid Card Number" (row 3) test
#     cases are corrected:
#       * Row2 H (ContainsValidation)    -> "balance"
#       * Row2 I (NotContainsValidation) -> 'errorCode": "PY-1501"'
#       * Row3 H (ContainsValidation)    -> mod10/invalid format JSON message
#   - Several header/body cells in columns B, C, H, I pick up the workbook's
#     "text" number format (same style already used by columns A and D-G),
#     and the H/I validation columns get word-wrap turned on so the longer
#     JSON snippets display legibly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GiftCards")

# --- Row 1 (header) -------------------------------------------------------
# B1, H1, I1 gain the shared "text" number format used elsewhere in the row.
$ws.Range("B1").NumberFormat = "@"
$ws.Range("H1").NumberFormat = "@"
$ws.Range("I1").NumberFormat = "@"

# --- Row 2 ------------------------------------------------------------------
# A2: TCNumber switches from a numeric 1 to the text value "1".
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("C2").NumberFormat = "@"

# H3 must be written before H2/I2 so that new shared-string entries line up
# in the same order the original workbook recorded them
# (idx49="1" from A2, idx50=mod10 JSON from H3, idx51="balance" from H2,
#  idx52='errorCode": "PY-1501"' from I2).
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").WrapText = $true
$ws.Range("H3").Value = "{`n  ""validationMessage"": ""The cardNumber specified did not pass mod10 check or is invalid format.""`n}"

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").WrapText = $true
$ws.Range("H2").Value = "balance"

$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "errorCode"": ""PY-1501"""

# --- Row 3 --------------------------------------------------------------
$ws.Range("B3").NumberFormat = "@"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").WrapText = $true

# --- Row 4 --------------------------------------------------------------
$ws.Range("B4").NumberFormat = "@"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").WrapText = $true

# --- Row 5 --------------------------------------------------------------
$ws.Range("B5").NumberFormat = "@"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").WrapText = $true

# --- Row 6 --------------------------------------------------------------
$ws.Range("B6").NumberFormat = "@"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").WrapText = $true

# --- Row heights: the sheet switches from an automatic ht="16" (rows 2-6)
# / unset (row 1) to an explicit custom height of 15 for every row.
$ws.Rows.Item(1).RowHeight = 15
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(3).RowHeight = 15
$ws.Rows.Item(4).RowHeight = 15
$ws.Rows.Item(5).RowHeight = 15
$ws.Rows.Item(6).RowHeight = 15

# --- Selection: the author left the cursor on T15 before saving. ----------
$ws.Activate() | Out-Null
$ws.Range("T15").Select() | Out-Null

Write-Output "GiftCards sheet updated"
